# Split the single "Lead Time (Days)" column into two columns:
#   F = "Lead Time Oct. (Days)"  (existing data, relabeled)
#   G = "Lead Time Sept. (Days)" (new column inserted after F)
# and populate the new column with data, plus add an averages row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current MTTR column (G), i.e. right after
# the existing "Lead Time (Days)" column (F). Everything from G onward
# shifts one column to the right (G->H, H->I, etc.), matching the diff.
$ws.Columns("G").Insert()

# Relabel the header of the (now-split) lead-time columns.
$ws.Range("F1").Value = "Lead Time Oct. (Days)"
$ws.Range("G1").Value = "Lead Time Sept. (Days)"

# Try to match the original column's custom width as closely as the engine
# allows (character-width grid rounding means this can't be bit-exact).
$ws.Columns("G").ColumnWidth = 17.96

# New "Lead Time Sept. (Days)" values for the four data rows.
$ws.Range("G2").Value = 38
$ws.Range("G3").Value = 75
$ws.Range("G4").Value = 73
$ws.Range("G5").Value = 25

# New averages row under the data, with an integer number format.
$ws.Range("F8").Formula = "=AVERAGE(F2:F7)"
$ws.Range("G8").Formula = "=AVERAGE(G2:G7)"
$ws.Range("F8:G8").NumberFormat = "0"

# Update the remembered selection to match the saved workbook state.
$ws.Range("F15").Select()
